{"js": "// Replace the date title and every \"a OP b=\" expression in the table\n// with the new values from the commit, matched positionally (each\n// paragraph in the document body holds exactly one such text run, in\n// document order: the title paragraph first, then every table cell).\nconst replacements = [\n  [\"2024-04-13 Saturday\", \"2024-04-14 Sunday\"],\n  [\"88+2=\", \"15+4=\"],\n  [\"7+24=\", \"6+28=\"],\n  [\"85-26=\", \"46+45=\"],\n  [\"80-68=\", \"18+1=\"],\n  [\"65-5=\", \"99-89=\"],\n  [\"49+47=\", \"11+48=\"],\n  [\"88-34=\", \"58-23=\"],\n  [\"14+35=\", \"31-10=\"],\n  [\"59-31=\", \"64-6=\"],\n  [\"29-21=\", \"6+57=\"],\n  [\"63-45=\", \"16+0=\"],\n  [\"0+73=\", \"85-14=\"],\n  [\"67-21=\", \"10+66=\"],\n  [\"49-45=\", \"59-8=\"],\n  [\"51-0=\", \"60+26=\"],\n  [\"48+45=\", \"54-18=\"],\n  [\"92-34=\", \"79-43=\"],\n  [\"47+32=\", \"57-9=\"],\n  [\"43-0=\", \"43+56=\"],\n  [\"18-4=\", \"39+56=\"],\n  [\"18+55=\", \"41+4=\"],\n  [\"15+47=\", \"64-29=\"],\n  [\"17+13=\", \"88+11=\"],\n  [\"67-57=\", \"24-23=\"],\n  [\"86-81=\", \"90-78=\"],\n  [\"81-17=\", \"46+33=\"],\n  [\"35+49=\", \"51+5=\"],\n  [\"55-13=\", \"98-22=\"],\n  [\"34+33=\", \"71-68=\"],\n  [\"28-3=\", \"16+73=\"],\n  [\"96-21=\", \"15+19=\"],\n  [\"22+70=\", \"95-87=\"],\n  [\"67-36=\", \"29-26=\"],\n  [\"98-30=\", \"26+35=\"],\n  [\"28+7=\", \"66-64=\"],\n  [\"96-8=\", \"16+67=\"],\n  [\"93-56=\", \"57-38=\"],\n  [\"88-27=\", \"9+66=\"],\n  [\"9-8=\", \"64-55=\"],\n  [\"22+57=\", \"82-76=\"],\n  [\"43+0=\", \"98-86=\"],\n  [\"11+86=\", \"6+91=\"],\n  [\"38+15=\", \"89-7=\"],\n  [\"47+50=\", \"72+11=\"],\n  [\"17+50=\", \"63-22=\"],\n  [\"71+0=\", \"17+12=\"],\n  [\"35+56=\", \"39+5=\"],\n  [\"31-21=\", \"74+17=\"],\n  [\"4+95=\", \"84-75=\"],\n  [\"86-16=\", \"69-40=\"],\n  [\"67-37=\", \"60-17=\"],\n  [\"27+2=\", \"80-41=\"],\n  [\"96-13=\", \"31+20=\"],\n  [\"23+33=\", \"58-50=\"],\n  [\"82-21=\", \"30-0=\"],\n  [\"80-36=\", \"2+42=\"],\n  [\"73-37=\", \"66-24=\"],\n  [\"47-18=\", \"5+78=\"],\n  [\"88-36=\", \"48-31=\"],\n  [\"72-25=\", \"62+29=\"],\n  [\"66-19=\", \"65-55=\"],\n  [\"56-31=\", \"16+74=\"],\n  [\"27+1=\", \"53-23=\"],\n  [\"90-88=\", \"72-1=\"],\n  [\"62-15=\", \"4+88=\"],\n  [\"20+76=\", \"13+41=\"],\n  [\"22+21=\", \"69+16=\"],\n  [\"59-10=\", \"53+35=\"],\n  [\"0+2=\", \"88-64=\"],\n  [\"20-4=\", \"16+83=\"],\n  [\"12-11=\", \"74-33=\"],\n  [\"2+74=\", \"69-27=\"],\n  [\"90-35=\", \"31+49=\"],\n  [\"72-20=\", \"43-28=\"],\n  [\"89+9=\", \"11+69=\"],\n  [\"91-44=\", \"71-9=\"],\n  [\"96-9=\", \"41-3=\"],\n  [\"77+12=\", \"31-18=\"],\n  [\"25+48=\", \"46-28=\"],\n  [\"79-10=\", \"60-38=\"],\n  [\"96-86=\", \"22+12=\"],\n  [\"69-33=\", \"62-17=\"],\n  [\"49-46=\", \"26+70=\"],\n  [\"26-22=\", \"15+28=\"],\n  [\"46-29=\", \"60-56=\"],\n  [\"85-69=\", \"75-11=\"],\n  [\"21+56=\", \"71+23=\"],\n  [\"18+16=\", \"45+40=\"],\n  [\"78+0=\", \"37+43=\"],\n  [\"56+33=\", \"73-25=\"],\n  [\"90-41=\", \"82+7=\"],\n  [\"36+41=\", \"27+6=\"],\n  [\"85-74=\", \"69+9=\"],\n  [\"74-63=\", \"26+72=\"],\n  [\"0+69=\", \"74-30=\"],\n  [\"49+27=\", \"16+44=\"],\n  [\"28+30=\", \"36+56=\"],\n  [\"16+41=\", \"99-1=\"],\n  [\"95-46=\", \"49-33=\"],\n  [\"97-42=\", \"87-19=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nif (paragraphs.items.length !== replacements.length) {\n  throw new Error(\n    \"Unexpected paragraph count: \" + paragraphs.items.length +\n    \" (expected \" + replacements.length + \")\"\n  );\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = replacements[i];\n  const actual = paragraphs.items[i].text;\n  if (actual !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \" + JSON.stringify(oldText) +\n      \" but found \" + JSON.stringify(actual)\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and every \"a OP b=\" expression inside the\n# 20x5 practice-problems table with the new values from the commit.\n# Each paragraph/cell is matched positionally against its current (old)\n# text before being overwritten, so the script fails loudly instead of\n# silently corrupting a cell if the document does not look as expected.\n# Word COM Range.Text includes the trailing paragraph mark (chr 13) and,\n# for table cells, the cell mark (chr 7) as well, so both are trimmed\n# off before comparing/using the text.\n\n$d = $word.ActiveDocument\n\nfunction Trim-RangeText($text) {\n    return $text.TrimEnd([char]13, [char]7)\n}\n\n$titleOld = '2024-04-13 Saturday'\n$titleNew = '2024-04-14 Sunday'\n\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = Trim-RangeText $titlePara.Range.Text\nif ($titleText -ne $titleOld) {\n    throw \"Title paragraph text mismatch: expected '$titleOld' but found '$titleText'\"\n}\n$titlePara.Range.Text = $titleNew\n\n$cellValues = @(\n    @('88+2=', '15+4='),\n    @('7+24=', '6+28='),\n    @('85-26=', '46+45='),\n    @('80-68=', '18+1='),\n    @('65-5=', '99-89='),\n    @('49+47=', '11+48='),\n    @('88-34=', '58-23='),\n    @('14+35=', '31-10='),\n    @('59-31=', '64-6='),\n    @('29-21=', '6+57='),\n    @('63-45=', '16+0='),\n    @('0+73=', '85-14='),\n    @('67-21=', '10+66='),\n    @('49-45=', '59-8='),\n    @('51-0=', '60+26='),\n    @('48+45=', '54-18='),\n    @('92-34=', '79-43='),\n    @('47+32=', '57-9='),\n    @('43-0=', '43+56='),\n    @('18-4=', '39+56='),\n    @('18+55=', '41+4='),\n    @('15+47=', '64-29='),\n    @('17+13=', '88+11='),\n    @('67-57=', '24-23='),\n    @('86-81=', '90-78='),\n    @('81-17=', '46+33='),\n    @('35+49=', '51+5='),\n    @('55-13=', '98-22='),\n    @('34+33=', '71-68='),\n    @('28-3=', '16+73='),\n    @('96-21=', '15+19='),\n    @('22+70=', '95-87='),\n    @('67-36=', '29-26='),\n    @('98-30=', '26+35='),\n    @('28+7=', '66-64='),\n    @('96-8=', '16+67='),\n    @('93-56=', '57-38='),\n    @('88-27=', '9+66='),\n    @('9-8=', '64-55='),\n    @('22+57=', '82-76='),\n    @('43+0=', '98-86='),\n    @('11+86=', '6+91='),\n    @('38+15=', '89-7='),\n    @('47+50=', '72+11='),\n    @('17+50=', '63-22='),\n    @('71+0=', '17+12='),\n    @('35+56=', '39+5='),\n    @('31-21=', '74+17='),\n    @('4+95=', '84-75='),\n    @('86-16=', '69-40='),\n    @('67-37=', '60-17='),\n    @('27+2=', '80-41='),\n    @('96-13=', '31+20='),\n    @('23+33=', '58-50='),\n    @('82-21=', '30-0='),\n    @('80-36=', '2+42='),\n    @('73-37=', '66-24='),\n    @('47-18=', '5+78='),\n    @('88-36=', '48-31='),\n    @('72-25=', '62+29='),\n    @('66-19=', '65-55='),\n    @('56-31=', '16+74='),\n    @('27+1=', '53-23='),\n    @('90-88=', '72-1='),\n    @('62-15=', '4+88='),\n    @('20+76=', '13+41='),\n    @('22+21=', '69+16='),\n    @('59-10=', '53+35='),\n    @('0+2=', '88-64='),\n    @('20-4=', '16+83='),\n    @('12-11=', '74-33='),\n    @('2+74=', '69-27='),\n    @('90-35=', '31+49='),\n    @('72-20=', '43-28='),\n    @('89+9=', '11+69='),\n    @('91-44=', '71-9='),\n    @('96-9=', '41-3='),\n    @('77+12=', '31-18='),\n    @('25+48=', '46-28='),\n    @('79-10=', '60-38='),\n    @('96-86=', '22+12='),\n    @('69-33=', '62-17='),\n    @('49-46=', '26+70='),\n    @('26-22=', '15+28='),\n    @('46-29=', '60-56='),\n    @('85-69=', '75-11='),\n    @('21+56=', '71+23='),\n    @('18+16=', '45+40='),\n    @('78+0=', '37+43='),\n    @('56+33=', '73-25='),\n    @('90-41=', '82+7='),\n    @('36+41=', '27+6='),\n    @('85-74=', '69+9='),\n    @('74-63=', '26+72='),\n    @('0+69=', '74-30='),\n    @('49+27=', '16+44='),\n    @('28+30=', '36+56='),\n    @('16+41=', '99-1='),\n    @('95-46=', '49-33='),\n    @('97-42=', '87-19=')\n)\n\n$table = $d.Tables.Item(1)\n$rows = $table.Rows.Count\n$cols = $table.Columns.Count\nif (($rows * $cols) -ne $cellValues.Count) {\n    throw \"Unexpected table size: $rows x $cols (expected $($cellValues.Count) cells)\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $pair = $cellValues[$idx]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $cell = $table.Cell($r, $c)\n        $cellText = Trim-RangeText $cell.Range.Text\n        if ($cellText -ne $oldText) {\n            throw \"Cell ($r,$c) text mismatch: expected '$oldText' but found '$cellText'\"\n        }\n        $cell.Range.Text = $newText\n        $idx++\n    }\n}\n"}
